# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 350 and 351) above the existing
# historical records for "Zapallo" / "Camote" at Vega Monumental
# Concepción, shifting the rest of the data block down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 350; existing rows 350:362 shift to 352:364
$ws.Rows("350:351").Insert()

# --- New row 350 ---
$ws.Cells.Item(350, 1).Value = 11
$ws.Cells.Item(350, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(350, 3).Value = "Bíobío"
$ws.Cells.Item(350, 4).Value = 45008
$ws.Cells.Item(350, 5).Value = 8
$ws.Cells.Item(350, 6).Value = 100112045
$ws.Cells.Item(350, 7).Value = "Zapallo"
$ws.Cells.Item(350, 8).Value = "Camote"
$ws.Cells.Item(350, 9).Value = "1a (cosecha)"
$ws.Cells.Item(350, 10).Value = 600
$ws.Cells.Item(350, 11).Value = 450
$ws.Cells.Item(350, 12).Value = 500
$ws.Cells.Item(350, 13).Value = 475
$ws.Cells.Item(350, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(350, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(350, 16).Value = 475
$ws.Cells.Item(350, 17).Value = 1
$ws.Cells.Item(350, 18).Value = "Hortaliza"

# --- New row 351 ---
$ws.Cells.Item(351, 1).Value = 11
$ws.Cells.Item(351, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(351, 3).Value = "Bíobío"
$ws.Cells.Item(351, 4).Value = 45008
$ws.Cells.Item(351, 5).Value = 8
$ws.Cells.Item(351, 6).Value = 100112045
$ws.Cells.Item(351, 7).Value = "Zapallo"
$ws.Cells.Item(351, 8).Value = "Camote"
$ws.Cells.Item(351, 9).Value = "2a (cosecha)"
$ws.Cells.Item(351, 10).Value = 300
$ws.Cells.Item(351, 11).Value = 350
$ws.Cells.Item(351, 12).Value = 350
$ws.Cells.Item(351, 13).Value = 350
$ws.Cells.Item(351, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(351, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(351, 16).Value = 350
$ws.Cells.Item(351, 17).Value = 1
$ws.Cells.Item(351, 18).Value = "Hortaliza"
